# Weekly refresh of fruit/vegetable data: re-order the date-keyed rows
# (2-13) by permuting the D/H/J/K/L/M/P values among rows, per the new
# source pull. Row 10 is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values (rows 2-13) for the columns that move.
$before = @{}
for ($r = 2; $r -le 13; $r++) {
    $before[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        H = $ws.Cells.Item($r, 8).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
}

# Mapping: new row -> old row whose values it now takes on.
$mapping = @{
    2 = 5
    3 = 8
    4 = 11
    5 = 2
    6 = 13
    7 = 9
    8 = 7
    9 = 12
    10 = 10
    11 = 3
    12 = 6
    13 = 4
}

foreach ($r in $mapping.Keys) {
    $src = $before[$mapping[$r]]
    $ws.Cells.Item($r, 4).Value2 = $src.D
    $ws.Cells.Item($r, 8).Value2 = $src.H
    $ws.Cells.Item($r, 10).Value2 = $src.J
    $ws.Cells.Item($r, 11).Value2 = $src.K
    $ws.Cells.Item($r, 12).Value2 = $src.L
    $ws.Cells.Item($r, 13).Value2 = $src.M
    $ws.Cells.Item($r, 16).Value2 = $src.P
}
